$wb = $excel.ActiveWorkbook

# 1. The "How is RTS Progresso blocked?" KPI row on the KPIs sheet currently has
#    Type = "Blocking All Shelves". That dedicated sheet/tab is being folded into
#    the regular "Blocking" sheet, so update the Type to "Blocking".
$kpis = $wb.Worksheets.Item("KPIs")
$kpis.Range("B14").Value = "Blocking"

# 2. Copy the single data row from "Blocking All Shelves" onto the end of the
#    "Blocking" sheet, marking its Result as "blocking covers" (covering all
#    shelves) instead of the old sheet's separate "blocking" result.
$blocking = $wb.Worksheets.Item("Blocking")
$blocking.Range("A8").Value = "How is RTS Progresso blocked?"
$blocking.Range("B8").Value = "GMI_SEGMENT"
$blocking.Range("C8").Value = "RTS"
$blocking.Range("D8").Value = "GMI_BRAND"
$blocking.Range("E8").Value = "PROGRESSO"
$blocking.Range("F8").Value = "blocking covers"

# 3. Remove the now-redundant "Blocking All Shelves" sheet entirely.
$allShelves = $wb.Worksheets.Item("Blocking All Shelves")
$allShelves.Delete()
